# Natmi following Dr Hou advice
# Recompute the Tnfsf11 -> Tnfrsf11a LR-pair table: update the existing
# FAPs-sourced rows (2-6) with refreshed statistics, and add five new rows
# (7-11) describing edges sourced from "sCs" (the sixth target cluster that
# used to be a target only, now also used as a sending cluster).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol |
# D Target cluster | E..T the various expression / specificity metrics.
$data = @(
    @("FAPs", "Tnfsf11", "Tnfrsf11a", "ECs",    3, 1, 1.409861,  4.229583,  0.8996799439973785, 0.9308061179577465, 2, 1, 2.216812,           4.433624,           0.03478955304943803, 0.0236002455591832,  3.125396783132,     18.752380698792,    0.03129946313921224,  0.02196725295179286),
    @("FAPs", "Tnfsf11", "Tnfrsf11a", "M1",     3, 1, 1.409861,  4.229583,  0.8996799439973785, 0.9308061179577465, 3, 1, 30.291975,          90.875925,          0.4753872999761597,  0.4837338812262644,  42.707474165475,    384.367267489275,   0.4276964194196163,   0.4502624561088528),
    @("FAPs", "Tnfsf11", "Tnfrsf11a", "M2",     3, 1, 1.409861,  4.229583,  0.8996799439973785, 0.9308061179577465, 3, 1, 28.90471933333333, 86.714158,          0.4536163944556945,  0.4615807344641337,  40.75163650401267,  366.764728536114,   0.408109572360192,    0.4296421715706457),
    @("FAPs", "Tnfsf11", "Tnfrsf11a", "Neutro", 3, 1, 1.409861,  4.229583,  0.8996799439973785, 0.9308061179577465, 3, 1, 1.225528333333333, 3.676585,           0.01923283659872347, 0.01957051586224036, 1.727824601561667,  15.550421414055,    0.01730339735405026,  0.01821635589616245),
    @("FAPs", "Tnfsf11", "Tnfrsf11a", "sCs",    3, 1, 1.409861,  4.229583,  0.8996799439973785, 0.9308061179577465, 2, 1, 1.0815885,          2.163177,           0.01697391591998424, 0.01151462288817844, 1.5248894441985,    9.149336665191001,  0.01527109172430763,  0.01071788143029279),
    @("sCs",  "Tnfsf11", "Tnfrsf11a", "ECs",    1, 0.5, 0.1572085, 0.314417, 0.1003200560026214, 0.06919388204225353, 2, 1, 2.216812,           4.433624,           0.03478955304943803, 0.0236002455591832,  0.348501689302,     1.394006757208,     0.003490089910225793, 0.00163299260739034),
    @("sCs",  "Tnfsf11", "Tnfrsf11a", "M1",     1, 0.5, 0.1572085, 0.314417, 0.1003200560026214, 0.06919388204225353, 3, 1, 30.291975,          90.875925,          0.4753872999761597,  0.4837338812262644,  4.762155951787499,  28.572935710725,    0.04769088055654334,  0.03347142511741162),
    @("sCs",  "Tnfsf11", "Tnfrsf11a", "M2",     1, 0.5, 0.1572085, 0.314417, 0.1003200560026214, 0.06919388204225353, 3, 1, 28.90471933333333, 86.714158,          0.4536163944556945,  0.4615807344641337,  4.544067569314334,  27.264405415886,    0.0455068220955025,   0.03193856289348802),
    @("sCs",  "Tnfsf11", "Tnfrsf11a", "Neutro", 1, 0.5, 0.1572085, 0.314417, 0.1003200560026214, 0.06919388204225353, 3, 1, 1.225528333333333, 3.676585,           0.01923283659872347, 0.01957051586224036, 0.1926634709908333, 1.155980825945,     0.001929439244673206, 0.001354159966077911),
    @("sCs",  "Tnfsf11", "Tnfrsf11a", "sCs",    1, 0.5, 0.1572085, 0.314417, 0.1003200560026214, 0.06919388204225353, 2, 1, 1.0815885,          2.163177,           0.01697391591998424, 0.01151462288817844, 0.17003490570225,   0.6801396228090001, 0.001702824195676607, 0.0007967414578856515)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $colNum = 1 + $j
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}

Write-Output "Wrote $($data.Length) data rows (rows 2-$(1 + $data.Length)) to $($ws.Name)"
